# Auto-generated Excel COM-interop script applying the Linea 141 schedule update
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 'Última actualización: 08:42:31'
$ws.Range("A3").Value = 'Total filas: 101'

$rows = @(
  @(79, '08:42:31', '08:51', '14_ABASTO', 9, 'LP1912'),
  @(80, '07:31:43', '08:52', '81_EL PELIGRO', 81, 'LP1912'),
  @(81, '08:42:31', '08:53', '10_OLMOS', 11, 'LP1912'),
  @(82, '06:59:44', '08:54', '17_ROMERO', 115, 'LP1912'),
  @(83, '07:31:43', '09:01', '215A_EL PATO', 90, 'LP1912'),
  @(84, '07:57:27', '09:03', '11_ETCHEVERRY', 66, 'LP1912'),
  @(85, '07:31:43', '09:10', '16_P MOR-SANTA ANA', 99, 'LP1912'),
  @(86, '08:42:31', '09:10', '16_SANTA ANA', 28, 'LP1912'),
  @(87, '08:42:31', '09:11', '23_HERNANDEZ', 29, 'LP1912'),
  @(88, '07:31:43', '09:16', '27_EL RETIRO', 105, 'LP1912'),
  @(89, '07:57:27', '09:17', '27_EL RETIRO', 80, 'LP1912'),
  @(90, '07:57:27', '09:18', '81_EL PELIGRO', 81, 'LP1912'),
  @(91, '07:31:43', '09:21', '26_HERNANDEZ', 110, 'LP1912'),
  @(92, '07:31:43', '09:22', '17_ROMERO', 111, 'LP1912'),
  @(93, '08:20:43', '09:22', '16_SANTA ANA', 62, 'LP1912'),
  @(94, '07:31:43', '09:23', '11_ETCHEVERRY', 112, 'LP1912'),
  @(95, '07:57:27', '09:23', '17_ROMERO', 86, 'LP1912'),
  @(96, '07:57:27', '09:32', '15_ABASTO', 95, 'LP1912'),
  @(97, '07:57:27', '09:33', '10_OLMOS', 96, 'LP1912'),
  @(98, '08:42:31', '09:34', '23_HERNANDEZ', 52, 'LP1912'),
  @(99, '08:20:43', '09:41', '215C_EL PATO', 81, 'LP1912'),
  @(100, '07:57:27', '09:42', '215C_EL PATO', 105, 'LP1912'),
  @(101, '07:57:27', '09:43', '14_ABASTO', 106, 'LP1912'),
  @(102, '08:42:31', '10:07', '10_OLMOS', 85, 'LP1912'),
  @(103, '08:20:43', '10:08', '10_OLMOS', 108, 'LP1912'),
  @(104, '08:20:43', '10:12', '15_ABASTO', 112, 'LP1912'),
  @(105, '08:42:31', '10:21', '26_HERNANDEZ', 99, 'LP1912'),
  @(106, '08:42:31', '10:26', '215A_EL PATO', 104, 'LP1912')
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)

$ws.Range("A2").Value = 'Última actualización: 08:42:31'
$ws.Range("A3").Value = 'Total filas: 14'

$rows = ,@(19, '08:42:31', '10:26', '215A_EL PATO', 104, 'LP1912')
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)

$ws.Range("A2").Value = 'Última actualización: 08:42:31'
$ws.Range("A3").Value = 'Total filas: 22'

$rows = ,@(27, '08:42:31', '10:03', '215B_LP-P MOR-40 Y 115', 81, 'L6173')
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
